$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new Price value would otherwise be auto-parsed as a number by
# Excel (losing formatting like trailing zeros); force them to stay text first.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '24.837.23'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '1.687.16'
$ws.Range('E3').Value = '  -1.20%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.92%  '
$ws.Range('D5').Value = '314.53'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('D7').Value = '0.3932'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.3972'
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('D9').Value = '1.005'
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('D10').Value = '1.425'
$ws.Range('E10').Value = '  -5.03%  '
$ws.Range('D11').Value = '51.90'
$ws.Range('E11').Value = '  -3.29%  '
$ws.Range('D12').Value = '0.08682'
$ws.Range('E12').Value = '  -1.68%  '
$ws.Range('D13').Value = '25.19'
$ws.Range('E13').Value = '  -4.21%  '
$ws.Range('D14').Value = '7.297'
$ws.Range('E14').Value = '  -2.71%  '
$ws.Range('D15').Value = '7.808'
$ws.Range('E15').Value = '  -4.05%  '
$ws.Range('D16').Value = '0.00001319'
$ws.Range('E16').Value = '  -3.30%  '
$ws.Range('D17').Value = '1.604.44'
$ws.Range('E17').Value = '  -5.96%  '
$ws.Range('D18').Value = '94.08'
$ws.Range('E18').Value = '  -3.79%  '
$ws.Range('D19').Value = '0.07122'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').Value = '20.21'
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').Value = '7.172'
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').Value = '14.13'
$ws.Range('E23').Value = '  -1.93%  '
$ws.Range('D24').Value = '24.845.01'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = '2.400'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').Value = '2.785'
$ws.Range('E26').Value = '  -8.21%  '
$ws.Range('D27').Value = '23.34'
$ws.Range('E27').Value = '  +1.22%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '161.60'
$ws.Range('E28').Value = '  -3.76%  '
$ws.Range('B29').Value = 'HuobiToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D29').Value = '5.881'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('D30').Value = '149.48'
$ws.Range('E30').Value = '  +2.94%  '
$ws.Range('D31').Value = '2.640'
$ws.Range('E31').Value = '  +21.19%  '
$ws.Range('D32').Value = '7.858'
$ws.Range('E32').Value = '  -8.65%  '
$ws.Range('D33').Value = '1.814.27'
$ws.Range('E33').Value = '  -3.93%  '
$ws.Range('D34').Value = '0.08448'
$ws.Range('E34').Value = '  -4.49%  '
$ws.Range('D35').Value = '0.03080'
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('E36').Value = '  -5.17%  '
$ws.Range('D37').Value = '6.939'
$ws.Range('E37').Value = '  -4.14%  '
$ws.Range('D38').Value = '0.2821'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').Value = '0.09555'
$ws.Range('E39').Value = '  +3.54%  '
$ws.Range('E40').Value = '  -3.10%  '
$ws.Range('D41').Value = '0.7993'
$ws.Range('E41').Value = '  -6.01%  '
$ws.Range('E42').Value = '  -3.34%  '
$ws.Range('D43').Value = '1.459'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('D44').Value = '16.76'
$ws.Range('E44').Value = '  -4.92%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.7191'
$ws.Range('E45').Value = '  -4.28%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '2.592'
$ws.Range('E46').Value = '  -4.40%  '
$ws.Range('D47').Value = '4.198'
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('D48').Value = '0.08728'
$ws.Range('E48').Value = '  +5.59%  '
$ws.Range('D49').Value = '1.004'
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('D50').Value = '1.343'
$ws.Range('E50').Value = '  -3.87%  '
$ws.Range('D51').Value = '138.88'
$ws.Range('E51').Value = '  -1.36%  '
